$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# LOQ4234.xlsx restructuring
#   * Row 10 (Objetivos:) gets a new Portuguese objectives paragraph in B/C.
#   * Starting at row 13, every label that used to sit in column A moves down
#     one row (A13->A14, A14->A15, ... A21->A22); three brand new Portuguese
#     bodies are introduced (Programa resumido, Programa, Bibliografia) and
#     the "Docentes responsaveis" value moves up into row 13's B/C (row 13
#     keeps no label at all). The three rows whose body text does not change
#     semantically (Metodo/Criterio/Norma de recuperacao content) simply stay
#     on the same row number while only their A-label text is swapped out.
# ---------------------------------------------------------------------------

# ---- 1. Capture every old value we still need, before anything is written --
$docentesValue      = $ws.Range("B18").Value2   # "5840560 - Marco Antonio Carvalho Pereira"

$progResumidoLabel  = $ws.Range("A13").Value2   # "Programa resumido:"
$shortSyllabusLabel = $ws.Range("A14").Value2   # "Short syllabus:"
$shortSyllabusBody  = $ws.Range("B14").Value2   # English short-syllabus paragraph

$programaLabel      = $ws.Range("A15").Value2   # "Programa:"
$syllabusLabel      = $ws.Range("A16").Value2   # "Syllabus:"
$syllabusBody       = $ws.Range("B16").Value2   # English full syllabus paragraph

$avaliacaoLabel     = $ws.Range("A17").Value2   # "Avaliação:"
$metodoLabel        = $ws.Range("A18").Value2   # "Método:"
$criterioLabel      = $ws.Range("A19").Value2   # "Critério:"
$normaLabel         = $ws.Range("A20").Value2   # "Norma de recuperação:"
$bibliografiaLabel  = $ws.Range("A21").Value2   # "Bibliografia:"

# ---- 2. New literal text introduced by this edit --------------------------
$objetivosPt = 'Fomentar a cultura do empreendedorismo; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup ao longo do semestre.'

$programaResumidoPt = 'Características do Comportamento Empreendedor; Modelo de Negócios; Produto mínimo viável; Plano de Negócios.'

$programaPt = '1.Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2.Modelo de Negócios (Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.3.Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente.4.Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.'

$bibliografiaBody = 'BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.'

# ---- 3. Row 10: new Objetivos paragraph (PT) -------------------------------
$ws.Range("B10").Value2 = $objetivosPt
$ws.Range("C10").Value2 = $objetivosPt

# ---- 4. Build the new row 22 (Bibliografia) from a copy of row 16's format
$ws.Range("A16:C16").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows(22).RowHeight = 120
$ws.Range("A22").Value2 = $bibliografiaLabel
$ws.Range("B22").Value2 = $bibliografiaBody
$ws.Range("C22").Value2 = $bibliografiaBody

# ---- 5. Row 21: Norma de recuperação: label, body stays, height 60 --------
$ws.Range("A21").Value2 = $normaLabel
$ws.Rows(21).RowHeight = 60

# ---- 6. Row 20: Critério: label, body stays -------------------------------
$ws.Range("A20").Value2 = $criterioLabel

# ---- 7. Row 19: Método: label, body stays ----------------------------------
$ws.Range("A19").Value2 = $metodoLabel

# ---- 8. Row 18: Avaliação: label only; clear B/C ---------------------------
$ws.Range("A18").Value2 = $avaliacaoLabel
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows(18).AutoFit()

# ---- 9. Row 17: Syllabus: label + English syllabus body, height 120 -------
$ws.Range("B16:C16").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A17").Value2 = $syllabusLabel
$ws.Range("B17").Value2 = $syllabusBody
$ws.Range("C17").Value2 = $syllabusBody
$ws.Rows(17).RowHeight = 120

# ---- 10. Row 16: Programa: label + new Portuguese full programa body ------
$ws.Range("A16").Value2 = $programaLabel
$ws.Range("B16").Value2 = $programaPt
$ws.Range("C16").Value2 = $programaPt

# ---- 11. Row 15: Short syllabus: label + English body, height 60 ----------
$ws.Range("A15").Value2 = $shortSyllabusLabel
$ws.Range("B15").Value2 = $shortSyllabusBody
$ws.Range("C15").Value2 = $shortSyllabusBody
$ws.Rows(15).RowHeight = 60

# ---- 12. Row 14: Programa resumido: label + new PT one-liner --------------
$ws.Range("A14").Value2 = $progResumidoLabel
$ws.Range("B14").Value2 = $programaResumidoPt
$ws.Range("C14").Value2 = $programaResumidoPt

# ---- 13. Row 13: no label any more; Docentes responsáveis value moves in --
$ws.Range("A13").Clear()
$ws.Range("B13").Value2 = $docentesValue
$ws.Range("C13").Value2 = $docentesValue
$ws.Rows(13).AutoFit()
